# Biggest Atlantic Hurricanes - add hurricane-category SWITCH formulas to column E
# and update the sheet's saved selection, matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Biggest Atlantic Hurricanes")

# E4:E24 each get a SWITCH() formula classifying the Saffir-Simpson hurricane
# category. Note the source data authored this by writing the formula in E4
# referencing A2 (two rows above its own row) and then filling down through
# E24, so every row's formula keeps that same two-row-up offset into column A.
for ($row = 4; $row -le 24; $row++) {
    $refRow = $row - 2
    $formula = "=_xlfn.SWITCH(TRUE,A$refRow>=157,""Category 5 - Catastrophic"",A$refRow>=130,""Category 4 - Extreme"",A$refRow>=111,""Category 3 - Devastating"",A$refRow>=96,""Category 2 - Very Dangerous"",A$refRow>=74,""Category 1 - Dangerous"",TRUE,""Not a Hurricane"")"
    $ws.Range("E$row").Formula = $formula
}

# These cells previously only carried a "horizontal right" alignment style
# with no value; clear that formatting back to the workbook's Normal style
# now that they hold real formula content.
$ws.Range("E4:E24").Style = "Normal"

# Restore the saved selection rectangle recorded in the sheet view.
[void]$ws.Range("A1:G34").Select()
